# Update the assumption inputs on the "Main" and "Model" sheets.
$wb = $excel.ActiveWorkbook

# Main sheet: N3 (implied multiple / units) 43 -> 50
$wsMain = $wb.Worksheets.Item("Main")
$wsMain.Range("N3").Value = 50

# Model sheet: AC55 (perpetuity growth rate) 0.02 -> 0.01
# Model sheet: V58 (first-period assumption) 0.125 -> 0.128
$wsModel = $wb.Worksheets.Item("Model")
$wsModel.Range("AC55").Value = 0.01
$wsModel.Range("V58").Value = 0.128

# Move the active selection on Model to match the author's last cursor position.
$wsModel.Activate()
$wsModel.Range("AA58").Select()
